$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("K64").Value = 4485.7144
$ws.Range("J64").Value = 59566.668
$ws.Range("L64").Value = 59566.668
$ws.Range("M64").Value = -4237.7144
$ws.Range("N64").Value = -60062.668
$ws.Range("I64").Value = 4485.7144
$ws.Range("H64").Value = 35468.75

$ws.Range("I67").Value = 4485.7144
$ws.Range("L67").Value = 59566.668
$ws.Range("H67").Value = 35468.75
$ws.Range("M67").Value = -3627.7144
$ws.Range("K67").Value = 4485.7144
$ws.Range("N67").Value = -61282.668
$ws.Range("J67").Value = 59566.668

$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 7966.6665
$ws.Range("K69").Value = 0
$ws.Range("H69").Value = 7966.6665
$ws.Range("N69").Value = -25647.9995
$ws.Range("M69").ClearContents()
$ws.Range("L69").Value = 23899.9995

$ws.Range("I72").Value = 0
$ws.Range("L72").Value = 71699.9985
$ws.Range("N72").Value = -80435.9985
$ws.Range("M72").ClearContents()
$ws.Range("H72").Value = 7966.6665
$ws.Range("K72").Value = 0
$ws.Range("J72").Value = 7966.6665

$ws.Range("K113").Value = 2875
$ws.Range("I113").Value = 2875
$ws.Range("L113").Value = 2135.3333
$ws.Range("N113").Value = -8643.3333
$ws.Range("M113").Value = 379
$ws.Range("J113").Value = 2135.3333
$ws.Range("H113").Value = 2431.2

$ws.Range("M115").Value = -291.125
$ws.Range("K115").Value = 1858.125
$ws.Range("L115").Value = 6012
$ws.Range("H115").Value = 773.2222
$ws.Range("I115").Value = 619.375
$ws.Range("J115").Value = 2004
$ws.Range("N115").Value = -9146

$ws.Range("J137").Value = 142861200
$ws.Range("I137").Value = 1233.0834
$ws.Range("N137").Value = -428588700
$ws.Range("L137").Value = 428583600
$ws.Range("M137").Value = -1149.2502
$ws.Range("H137").Value = 32259936
$ws.Range("K137").Value = 3699.2502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M4").Value = -204.25
$ws.Range("N4").Value = -662
$ws.Range("I4").Value = 320.25
$ws.Range("L4").Value = 430
$ws.Range("K4").Value = 320.25
$ws.Range("J4").Value = 430
$ws.Range("H4").Value = 375.125

$ws.Range("N7").Value = -27628.5
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 27400.5
$ws.Range("J7").Value = 27400.5
$ws.Range("M7").ClearContents()
$ws.Range("I7").Value = 0
$ws.Range("H7").Value = 27400.5

$ws.Range("M32").Value = -6329.2036
$ws.Range("H32").Value = 8771.5
$ws.Range("I32").Value = 6616.2036
$ws.Range("K32").Value = 6616.2036

$ws.Range("H61").Value = 2389.7693
$ws.Range("M61").Value = -1174.7142
$ws.Range("N61").Value = -7026.6
$ws.Range("J61").Value = 6602.6
$ws.Range("I61").Value = 1386.7142
$ws.Range("L61").Value = 6602.6
$ws.Range("K61").Value = 1386.7142

$ws.Range("H102").Value = 1301.25
$ws.Range("M102").Value = 320.75
$ws.Range("I102").Value = 1301.25
$ws.Range("K102").Value = 1301.25

$ws.Range("K110").Value = 839.86957
$ws.Range("H110").Value = 951.36365
$ws.Range("I110").Value = 839.86957
$ws.Range("M110").Value = 1205.13043
$ws.Range("J110").Value = 1207.8
$ws.Range("L110").Value = 1207.8
$ws.Range("N110").Value = -5297.8

$ws.Range("K132").Value = 4990.200000000001
$ws.Range("J132").Value = 2041.8334
$ws.Range("I132").Value = 1663.4
$ws.Range("H132").Value = 1771.5238
$ws.Range("L132").Value = 6125.5002
$ws.Range("N132").Value = -11185.5002
$ws.Range("M132").Value = -2460.200000000001

$ws.Range("K136").Value = 4160.142599999999
$ws.Range("N136").Value = -24907.8
$ws.Range("J136").Value = 6602.6
$ws.Range("M136").Value = -1610.142599999999
$ws.Range("L136").Value = 19807.8
$ws.Range("H136").Value = 2389.7693
$ws.Range("I136").Value = 1386.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("L94").Value = 1025
$ws.Range("H94").Value = 513.6111
$ws.Range("I94").Value = 367.5
$ws.Range("N94").Value = -1927
$ws.Range("J94").Value = 1025
$ws.Range("M94").Value = 83.5
$ws.Range("K94").Value = 367.5

$ws.Range("N105").Value = -4994
$ws.Range("L105").Value = 1500
$ws.Range("H105").Value = 4042.8572
$ws.Range("J105").Value = 1500
$ws.Range("K105").Value = 5060
$ws.Range("M105").Value = -3313
$ws.Range("I105").Value = 5060

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J31").Value = 1469
$ws.Range("L31").Value = 1469
$ws.Range("K31").Value = 1259.6842
$ws.Range("M31").Value = -964.6841999999999
$ws.Range("I31").Value = 1259.6842
$ws.Range("H31").Value = 1309.92
$ws.Range("N31").Value = -2059

$ws.Range("I34").Value = 1259.6842
$ws.Range("J34").Value = 1469
$ws.Range("M34").Value = -1057.6842
$ws.Range("H34").Value = 1309.92
$ws.Range("K34").Value = 1259.6842
$ws.Range("L34").Value = 1469
$ws.Range("N34").Value = -1873

$ws.Range("K132").Value = 5871
$ws.Range("I132").Value = 1957
$ws.Range("H132").Value = 2537.0833
$ws.Range("M132").Value = -3341

$ws.Range("H134").Value = 4716
$ws.Range("J134").Value = 2750
$ws.Range("K134").Value = 16507.2
$ws.Range("M134").Value = -13972.2
$ws.Range("I134").Value = 5502.4
$ws.Range("N134").Value = -13320
$ws.Range("L134").Value = 8250

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("N105").Value = -2724005242
$ws.Range("L105").Value = 2724000000
$ws.Range("H105").Value = 908000000
$ws.Range("J105").Value = 908000000

$ws.Range("J107").Value = 43942.78
$ws.Range("N107").Value = -135668.34
$ws.Range("H107").Value = 71857.75
$ws.Range("L107").Value = 131828.34

$ws.Range("H131").Value = 1963511.9
$ws.Range("K131").Value = 16149
$ws.Range("I131").Value = 5383
$ws.Range("L131").Value = 7323312.600000001
$ws.Range("J131").Value = 2441104.2
$ws.Range("M131").Value = -11109
$ws.Range("N131").Value = -7333392.600000001

$ws.Range("J133").Value = 8128.5186
$ws.Range("H133").Value = 7317.879
$ws.Range("L133").Value = 24385.5558
$ws.Range("N133").Value = -34505.5558

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2005.4286
$ws.Range("M61").Value = -1025.6
$ws.Range("I61").Value = 1227.6
$ws.Range("K61").Value = 1227.6

$ws.Range("K113").Value = 1227.6
$ws.Range("I113").Value = 1227.6
$ws.Range("M113").Value = 942.4000000000001
$ws.Range("H113").Value = 2005.4286

$ws.Range("K132").Value = 4684.857
$ws.Range("J132").Value = 3583.1667
$ws.Range("I132").Value = 1561.619
$ws.Range("H132").Value = 2296.7273
$ws.Range("L132").Value = 10749.5001
$ws.Range("N132").Value = -15809.5001
$ws.Range("M132").Value = -2154.857

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("J20").Value = 29999
$ws.Range("I20").Value = 32499
$ws.Range("L20").Value = 29999
$ws.Range("N20").Value = -30479
$ws.Range("H20").Value = 32412.793
$ws.Range("K20").Value = 32499
$ws.Range("M20").Value = -32259

$ws.Range("K113").Value = 1213.23531
$ws.Range("I113").Value = 404.41177
$ws.Range("L113").Value = 1064.33331
$ws.Range("N113").Value = -5404.33331
$ws.Range("M113").Value = 956.76469
$ws.Range("J113").Value = 354.77777
$ws.Range("H113").Value = 387.23077

$ws.Range("K132").Value = 4439.4
$ws.Range("J132").Value = 3749.5
$ws.Range("I132").Value = 1479.8
$ws.Range("H132").Value = 2128.2856
$ws.Range("L132").Value = 11248.5
$ws.Range("N132").Value = -16308.5
$ws.Range("M132").Value = -1909.4
